# Refresh cryptocurrency price / volume snapshot pulled in by the
# scheduled GitHub Actions job. Column D (Price) holds locale-formatted
# text (dotted thousands separators), so force text formatting before
# writing to keep Excel from re-parsing values like "7.98" as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '66.168.29'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.09%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.557.08'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +4.53%  '

# Row 4
$ws.Range('E4').Value = '  -0.01%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '606.76'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.74%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.64'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.25%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.555.66'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +4.48%  '

# Row 8
$ws.Range('E8').Value = '  +0.09%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.496'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +5.87%  '

# Row 10
$ws.Range('E10').Value = '  +2.42%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.98'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.77%  '

# Row 12
$ws.Range('E12').Value = '  +2.39%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.166.39'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +4.71%  '

# Row 14
$ws.Range('E14').Value = '  +3.94%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '30.09'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.52%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.566.07'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +5.07%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '66.282.82'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.15%  '

# Row 18
$ws.Range('E18').Value = '  -0.47%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.36'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +10.30%  '

# Row 20
$ws.Range('E20').Value = '  +1.59%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.90'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.69%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '431.09'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +4.50%  '

# Row 23
$ws.Range('E23').Value = '  +6.34%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '79.13'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.94%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.700.57'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +4.54%  '

# Row 26
$ws.Range('E26').Value = '  -0.03%  '

# Row 27
$ws.Range('E27').Value = '  +9.25%  '

# Row 28
$ws.Range('E28').Value = '  +4.89%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.00'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +3.23%  '

# Row 30
$ws.Range('E30').Value = '  -0.78%  '

# Row 31
$ws.Range('E31').Value = '  +0.13%  '

# Row 32
$ws.Range('E32').Value = '  +2.10%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '25.52'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.73%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.554.37'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +4.57%  '

# Row 35
$ws.Range('E35').Value = '  -3.60%  '

# Row 36
$ws.Range('E36').Value = '  +0.06%  '

# Row 37
$ws.Range('E37').Value = '  +4.43%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '7.90'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +5.65%  '

# Row 39
$ws.Range('E39').Value = '  +2.15%  '

# Row 40
$ws.Range('E40').Value = '  +0.01%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '173.41'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.71%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0853'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.69%  '

# Row 43
$ws.Range('E43').Value = '  +4.04%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.894'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.82%  '

# Row 45
$ws.Range('E45').Value = '  +1.08%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '46.12'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.75%  '

# Row 47
$ws.Range('E47').Value = '  +2.30%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '25.78'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.03%  '

# Row 49
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '23.50'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +16.59%  '

# Row 50
$ws.Range('B50').Value = 'Cosmos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.14'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.69%  '

# Row 51
$ws.Range('B51').Value = 'dogwifhat'
$ws.Range('C51').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.35'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +4.79%  '
